$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has strain names in C2:C21 and, off to the side,
# a "Host Adaptation" column header in F1 with values in F2:F21.
# The author consolidated this into a single contiguous block in C:D,
# added a "Strain" label for the first column, and turned the range
# into an Excel Table (which supplies its own generic "Column1"/"Column2"
# header row above the existing labels).

# 1. Move the Host Adaptation column (F) next to the strain names (D).
$ws.Range("F1:F21").Cut($ws.Range("D1:D21"))

# 2. Insert a blank row above everything so there's room for the table's
#    auto-generated header row.
$ws.Range("C1:D1").Insert()

# 3. Label the (now second-row) strain column to match its neighbour's
#    existing "Host Adaptation" header.
$ws.Range("C2").Value = "Strain"

# 4. Select the data block and turn it into a Table.
$rng = $ws.Range("C1:D22")
$rng.Select()
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table4"

# 5. Match the column widths Excel applied to the new table columns.
$ws.Range("C1:D1").ColumnWidth = 10.14
